$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.538445
$ws.Range("H2").Value = 1.615335
$ws.Range("I2").Value = 0.03371608002174246
$ws.Range("J2").Value = 0.03371608002174246
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 1.091490115656667
$ws.Range("R2").Value = 9.823411040909999
$ws.Range("S2").Value = 0.0002224008582538388
$ws.Range("T2").Value = 0.0002224008582538387
$ws.Range("G3").Value = 0.538445
$ws.Range("H3").Value = 1.615335
$ws.Range("I3").Value = 0.03371608002174246
$ws.Range("J3").Value = 0.03371608002174246
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 138.0811539852066
$ws.Range("R3").Value = 1242.73038586686
$ws.Range("S3").Value = 0.02813526821222282
$ws.Range("T3").Value = 0.02813526821222282
$ws.Range("G4").Value = 0.538445
$ws.Range("H4").Value = 1.615335
$ws.Range("I4").Value = 0.03371608002174246
$ws.Range("J4").Value = 0.03371608002174246
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 26.29779684688834
$ws.Range("R4").Value = 236.680171621995
$ws.Range("S4").Value = 0.005358410951265804
$ws.Range("T4").Value = 0.005358410951265802
$ws.Range("I5").Value = 0.7539416098905094
$ws.Range("J5").Value = 0.7539416098905093
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 24.40733959722156
$ws.Range("R5").Value = 219.666056374994
$ws.Range("S5").Value = 0.00497321340454763
$ws.Range("T5").Value = 0.004973213404547629
$ws.Range("I6").Value = 0.7539416098905094
$ws.Range("J6").Value = 0.7539416098905093
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("Q6").Value = 3087.699621783192
$ws.Range("R6").Value = 27789.29659604873
$ws.Range("S6").Value = 0.6291463716109749
$ws.Range("T6").Value = 0.6291463716109748
$ws.Range("I7").Value = 0.7539416098905094
$ws.Range("J7").Value = 0.7539416098905093
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 588.0577836608149
$ws.Range("R7").Value = 5292.520052947334
$ws.Range("S7").Value = 0.1198220248749869
$ws.Range("T7").Value = 0.1198220248749869
$ws.Range("G8").Value = 3.391101666666666
$ws.Range("H8").Value = 10.173305
$ws.Range("I8").Value = 0.2123423100877482
$ws.Range("J8").Value = 0.2123423100877481
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 6.874154185392221
$ws.Range("R8").Value = 61.86738766852999
$ws.Range("S8").Value = 0.001400670302617147
$ws.Range("T8").Value = 0.001400670302617147
$ws.Range("G9").Value = 3.391101666666666
$ws.Range("H9").Value = 10.173305
$ws.Range("I9").Value = 0.2123423100877482
$ws.Range("J9").Value = 0.2123423100877481
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 869.628711223042
$ws.Range("R9").Value = 7826.658401007379
$ws.Range("S9").Value = 0.1771946158411398
$ws.Range("T9").Value = 0.1771946158411397
$ws.Range("G10").Value = 3.391101666666666
$ws.Range("H10").Value = 10.173305
$ws.Range("I10").Value = 0.2123423100877482
$ws.Range("J10").Value = 0.2123423100877481
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 165.6223063026761
$ws.Range("R10").Value = 1490.600756724085
$ws.Range("S10").Value = 0.03374702394399128
$ws.Range("T10").Value = 0.03374702394399127
